# Merge the three runs of the "g(n) - cost from start to current node"
# paragraph on slide 12 into a single run, as described by the diff.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(12)
$shp = $s.Shapes.Item(1)
$tf = $shp.TextFrame
$tr = $tf.TextRange

# Paragraph 11 (1-indexed) is the "g(n) ... to current node" bullet.
$para = $tr.Paragraphs(11, 1)

# Grab the paragraph's full character range and set its Text in one shot;
# this collapses the existing multiple runs into a single run while
# preserving the run-level formatting (rPr) of the paragraph.
$fullRange = $tr.Characters($para.Start, $para.Length)
$fullRange.Text = "g(n) " + [char]0x2013 + " cost from start to current node"
